# Siigo config: repoint the Config sheet from the old "escribana" sample
# credentials to the Siigo auth-token endpoint, add a header label, clear
# the sample secrets but keep the Input-style formatting, and leave a single
# hyperlink behind (on the apiUrl cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
# New header label in row 1.
$ws.Range("C1").Value = "Siigo config"

# Remove every existing hyperlink on the sheet (the old mailto: on C2 and
# the old auth-token link on C5) so we can re-create only the one we want.
$ws.Hyperlinks.Delete()

# C2 (username), C3 (accessKey), C4 (partnerId) lose their sample values but
# keep whatever formatting is already on them.
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# C5 (apiUrl) now holds the Siigo auth endpoint.
$ws.Range("C5").Value = "https://api.siigo.com/auth/token"

# --- formatting ---------------------------------------------------------
# Apply Excel's built-in "Input" cell style (orange fill + grey border) to
# the whole value column instead of the old "Hyperlink" look.
$ws.Range("C2:C5").Style = "Input"

# --- hyperlink ------------------------------------------------------------
# Re-add a single hyperlink, on the apiUrl cell, pointing at the same URL.
$ws.Hyperlinks.Add($ws.Range("C5"), "https://api.siigo.com/auth/token")

# --- view / sheet chrome -------------------------------------------------
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("B2").Select()
